{"js": "// Replace the five question headings per the diff:\n//   Q1: \"1. What is 5 + 7?\"                                  -> \"1. What is the area of a square with side 4?\"\n//   Q4: \"4. What is 20% of 50?\"                               -> \"4. What is 5 + 7?\"\n//   Q5: \"5. What is the area of a square with side 4?\"        -> \"5. What is 20% of 50?\"\n//   Q7: \"7. Solve: x/5 = 3\"                                   -> \"7. Solve: 2x + 3 = 11\"\n//   Q8: \"8. What is the perimeter of a square with side 6?\"   -> \"8. Solve: x/5 = 3\"\nconst replacements = [\n  [\"1. What is 5 + 7?\", \"1. What is the area of a square with side 4?\"],\n  [\"4. What is 20% of 50?\", \"4. What is 5 + 7?\"],\n  [\"5. What is the area of a square with side 4?\", \"5. What is 20% of 50?\"],\n  [\"7. Solve: x/5 = 3\", \"7. Solve: 2x + 3 = 11\"],\n  [\"8. What is the perimeter of a square with side 6?\", \"8. Solve: x/5 = 3\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Replace-QuestionText($oldText, $newText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n\n# 1. \"What is 5 + 7?\" moves to question 4; question 1 becomes the old question 5 text.\nReplace-QuestionText \"1. What is 5 + 7?\" \"1. What is the area of a square with side 4?\"\n\n# 4 <-> 5 swap of question text.\nReplace-QuestionText \"4. What is 20% of 50?\" \"4. What is 5 + 7?\"\nReplace-QuestionText \"5. What is the area of a square with side 4?\" \"5. What is 20% of 50?\"\n\n# 7 gets a brand-new equation; old question 8 text is replaced by old question 7 text.\nReplace-QuestionText \"7. Solve: x/5 = 3\" \"7. Solve: 2x + 3 = 11\"\nReplace-QuestionText \"8. What is the perimeter of a square with side 6?\" \"8. Solve: x/5 = 3\"\n"}
